$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 236, shifting existing rows 236-262 down to 237-263
$ws.Rows.Item(236).Insert()

# Fill in the newly inserted row 236 with the new weekly record
$ws.Range("A236").Value = 4
$ws.Range("B236").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C236").Value = "Los Lagos"
$ws.Range("D236").Value = 44946
$ws.Range("E236").Value = 10
$ws.Range("F236").Value = 100112009
$ws.Range("G236").Value = "Acelga"
$ws.Range("H236").Value = "Sin especificar"
$ws.Range("I236").Value = "Primera"
$ws.Range("J236").Value = 50
$ws.Range("K236").Value = 10000
$ws.Range("L236").Value = 10000
$ws.Range("M236").Value = 10000
$ws.Range("N236").Value = "$/docena de atados (12 kilos)"
$ws.Range("O236").Value = "Región de La Araucanía"
$ws.Range("P236").Value = 833
$ws.Range("Q236").Value = 12
$ws.Range("R236").Value = "Hortaliza"
